$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from column H to column I for each data row, then set values
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 2023
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = 5.7627222366917641
$ws.Range("H7").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = 12.150374768642443
$ws.Range("H8").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("I8").Value = 0.33526865592135835
$ws.Range("H9").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("H10").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 4.1862505557986136
$ws.Range("H11").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Value = 6.7483562655646434
$ws.Range("H12").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("H13").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I13").Value = 11.721778533441505
$ws.Range("H14").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = 0.39226026012037718
$ws.Range("H15").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("I15").Value = 15.431347214780089
$ws.Range("H16").Copy()
$ws.Range("I16").PasteSpecial(-4122)
$ws.Range("I16").Value = 12.56881331951053
$ws.Range("H17").Copy()
$ws.Range("I17").PasteSpecial(-4122)
$ws.Range("I17").Value = 4.9361801817513591
$ws.Range("H18").Copy()
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("I18").Value = 13.92788271827051
$ws.Range("H19").Copy()
$ws.Range("I19").PasteSpecial(-4122)
$ws.Range("I19").Value = 5.796504268446359
$ws.Range("H20").Copy()
$ws.Range("I20").PasteSpecial(-4122)
$ws.Range("I20").Value = 3.6469692666385813
$ws.Range("H21").Copy()
$ws.Range("I21").PasteSpecial(-4122)
$ws.Range("I21").Value = 1.2344990530700553
$ws.Range("H22").Copy()
$ws.Range("I22").PasteSpecial(-4122)
$ws.Range("H23").Copy()
$ws.Range("I23").PasteSpecial(-4122)
$ws.Range("I23").Value = 0.33854574252686492
$ws.Range("H24").Copy()
$ws.Range("I24").PasteSpecial(-4122)
$ws.Range("I24").Value = 7.0095823182535142
$ws.Range("H25").Copy()
$ws.Range("I25").PasteSpecial(-4122)
$ws.Range("I25").Value = 7.4180588363268161
$ws.Range("H26").Copy()
$ws.Range("I26").PasteSpecial(-4122)
$ws.Range("H27").Copy()
$ws.Range("I27").PasteSpecial(-4122)
$ws.Range("I27").Value = 1.3575537444685963
$ws.Range("H28").Copy()
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("I28").Value = 3.7000582818073822
$ws.Range("H29").Copy()
$ws.Range("I29").PasteSpecial(-4122)
$ws.Range("I29").Value = 7.0145832826742662
$ws.Range("H30").Copy()
$ws.Range("I30").PasteSpecial(-4122)
$ws.Range("I30").Value = 7.601388319014589
$ws.Range("H31").Copy()
$ws.Range("I31").PasteSpecial(-4122)
$ws.Range("I31").Value = 3.2001807961995414
$ws.Range("H32").Copy()
$ws.Range("I32").PasteSpecial(-4122)
$ws.Range("H33").Copy()
$ws.Range("I33").PasteSpecial(-4122)
$ws.Range("I33").Value = 6.1374688939827911
$ws.Range("H34").Copy()
$ws.Range("I34").PasteSpecial(-4122)
$ws.Range("I34").Value = 5.522716841454633
$ws.Range("H35").Copy()
$ws.Range("I35").PasteSpecial(-4122)
$ws.Range("I35").Value = 5.959494359842247
$ws.Range("H36").Copy()
$ws.Range("I36").PasteSpecial(-4122)
$ws.Range("I36").Value = 5.4831892692336535
$ws.Range("H37").Copy()
$ws.Range("I37").PasteSpecial(-4122)
$ws.Range("I37").Value = 5.7612749525079918
$excel.CutCopyMode = 0

# Remove stale selection on sheet view, select A1
$ws.Range("A1").Select()
